$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price (D) cells that would otherwise be
# auto-converted to numeric values by Excel, so the stored text matches exactly
# (preserves trailing zeros / literal formatting as in the source diff).
$ws.Range("D4:D9").NumberFormat = "@"
$ws.Range("D11:D17").NumberFormat = "@"
$ws.Range("D19:D20").NumberFormat = "@"
$ws.Range("D22:D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.014.60"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "1.744.13"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "250.31"
$ws.Range("E5").Value = "  +7.35%  "

$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").Value = "0.5166"
$ws.Range("E7").Value = "  -1.78%  "

$ws.Range("D8").Value = "0.2763"
$ws.Range("E8").Value = "  -0.90%  "

$ws.Range("D9").Value = "0.06200"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "1.744.37"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").Value = "0.07216"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "15.23"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "0.6512"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").Value = "4.642"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "77.98"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "26.046.28"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").Value = "11.86"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").Value = "0.000006814"
$ws.Range("E20").Value = "  +1.88%  "

$ws.Range("D21").Value = "1.966.42"
$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("D22").Value = "4.293"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "8.685"
$ws.Range("E23").Value = "  -1.50%  "

$ws.Range("D24").Value = "5.368"
$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("D25").Value = "135.87"
$ws.Range("E25").Value = "  -2.29%  "

$ws.Range("D26").Value = "1.513"
$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").Value = "15.29"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").Value = "1.784"
$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").Value = "105.76"
$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("D30").Value = "3.962"
$ws.Range("E30").Value = "  +4.99%  "

$ws.Range("D31").Value = "0.08288"
$ws.Range("E31").Value = "  -0.84%  "

$ws.Range("D33").Value = "0.04683"
$ws.Range("E33").Value = "  +3.08%  "

$ws.Range("D35").Value = "1.003"
$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").Value = "0.6253"
$ws.Range("E36").Value = "  -1.32%  "

$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  +1.35%  "

$ws.Range("D38").Value = "0.01606"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").Value = "1.939"
$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "0.9986"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").Value = "100.63"
$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("D42").Value = "0.3882"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").Value = "0.7565"
$ws.Range("E43").Value = "  +2.79%  "

$ws.Range("D44").Value = "5.025"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D47").Value = "55.58"
$ws.Range("E47").Value = "  +2.79%  "

$ws.Range("D48").Value = "0.05225"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("D49").Value = "30.70"
$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D50").Value = "7.629"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").Value = "0.3447"
$ws.Range("E51").Value = "  -0.58%  "

$ws.Range("E32").Value = "  +1.76%  "

$ws.Range("E34").Value = "  +0.93%  "

# Row 45/46: Aptos and Algorand swap positions with updated values
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "6.350"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1138"
$ws.Range("E46").Value = "  +0.31%  "
